# The document carries the Pearson/BTec logos twice each - once in the
# "default" header/footer and once in the "first page" header/footer.
# Each inline picture's name was swapped with its sibling's name:
#   Pearson logo (footers) : image1.png -> image2.png
#   BTec logo   (headers)  : image2.jpg -> image1.jpg
#
# Renaming an InlineShape only needs InlineShape.Name = "...". Selecting
# the shape first and then renaming it through Selection.InlineShapes is
# the reliable path for shapes that live in a footer story, so that
# pattern is used uniformly for headers and footers below.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlineLogo($range, $newName) {
    $shape = $range.InlineShapes.Item(1)
    $shape.Select()
    $word.Selection.InlineShapes.Item(1).Name = $newName
}

# Headers (BTec logo): image2.jpg -> image1.jpg
$headers = $sec.Headers
for ($i = 1; $i -le $headers.Count; $i++) {
    $h = $headers.Item($i)
    if ($h.Exists -and $h.Range.InlineShapes.Count -gt 0) {
        Rename-InlineLogo $h.Range "image1.jpg"
    }
}

# Footers (Pearson logo): image1.png -> image2.png
$footers = $sec.Footers
for ($i = 1; $i -le $footers.Count; $i++) {
    $f = $footers.Item($i)
    if ($f.Exists -and $f.Range.InlineShapes.Count -gt 0) {
        Rename-InlineLogo $f.Range "image2.png"
    }
}
